$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -20.40935985752628
$ws.Range("C2").Value = 1.8696478840477
$ws.Range("D2").Value = -20.40935985752628
$ws.Range("E2").Value = -20.40935985752628
$ws.Range("F2").Value = -20.40935985752628
$ws.Range("G2").Value = -20.40935985752628
$ws.Range("H2").Value = -20.40935985752628
$ws.Range("I2").Value = -20.40935985752628
$ws.Range("J2").Value = -20.40935985752628
$ws.Range("K2").Value = -20.40935985752628
$ws.Range("B3").Value = -20.40935985752628
$ws.Range("C3").Value = -20.40935985752628
$ws.Range("D3").Value = -20.40935985752628
$ws.Range("E3").Value = -20.40935985752628
$ws.Range("F3").Value = -20.40935985752628
$ws.Range("G3").Value = -20.40935985752628
$ws.Range("H3").Value = -20.40935985752628
$ws.Range("I3").Value = 1.276968396205431
$ws.Range("J3").Value = -20.40935985752628
$ws.Range("K3").Value = -20.40935985752628
$ws.Range("B4").Value = -20.40935985752628
$ws.Range("C4").Value = 1.978415847436124
$ws.Range("D4").Value = 1.692518701144461
$ws.Range("E4").Value = -20.40935985752628
$ws.Range("F4").Value = -20.40935985752628
$ws.Range("G4").Value = -20.40935985752628
$ws.Range("H4").Value = 1.482683723914448
$ws.Range("I4").Value = -20.40935985752628
$ws.Range("J4").Value = 1.058799288846011
$ws.Range("K4").Value = -20.40935985752628
$ws.Range("B5").Value = -20.40935985752628
$ws.Range("C5").Value = 1.702452388664
$ws.Range("D5").Value = -20.40935985752628
$ws.Range("E5").Value = -20.40935985752628
$ws.Range("F5").Value = -20.40935985752628
$ws.Range("G5").Value = 2.902442356457053
$ws.Range("H5").Value = -20.40935985752628
$ws.Range("I5").Value = -20.40935985752628
$ws.Range("J5").Value = -20.40935985752628
$ws.Range("K5").Value = -20.40935985752628
$ws.Range("B6").Value = -20.40935985752628
$ws.Range("C6").Value = -20.40935985752628
$ws.Range("D6").Value = -20.40935985752628
$ws.Range("E6").Value = -20.40935985752628
$ws.Range("F6").Value = -20.40935985752628
$ws.Range("G6").Value = -20.40935985752628
$ws.Range("H6").Value = -20.40935985752628
$ws.Range("I6").Value = -20.40935985752628
$ws.Range("J6").Value = -20.40935985752628
$ws.Range("K6").Value = -20.40935985752628
$ws.Range("B7").Value = 2.451037781537455
$ws.Range("C7").Value = -20.40935985752628
$ws.Range("D7").Value = -20.40935985752628
$ws.Range("E7").Value = -20.40935985752628
$ws.Range("F7").Value = -20.40935985752628
$ws.Range("G7").Value = -20.40935985752628
$ws.Range("H7").Value = -20.40935985752628
$ws.Range("I7").Value = -20.40935985752628
$ws.Range("J7").Value = -20.40935985752628
$ws.Range("K7").Value = -20.40935985752628
$ws.Range("B8").Value = -20.40935985752628
$ws.Range("C8").Value = -20.40935985752628
$ws.Range("D8").Value = -20.40935985752628
$ws.Range("E8").Value = 1.76356110268125
$ws.Range("F8").Value = -20.40935985752628
$ws.Range("G8").Value = -20.40935985752628
$ws.Range("H8").Value = -20.40935985752628
$ws.Range("I8").Value = -20.40935985752628
$ws.Range("J8").Value = -20.40935985752628
$ws.Range("K8").Value = -20.40935985752628
$ws.Range("B9").Value = 3.861150849231585
$ws.Range("C9").Value = -20.40935985752628
$ws.Range("D9").Value = -20.40935985752628
$ws.Range("E9").Value = -20.40935985752628
$ws.Range("F9").Value = -20.40935985752628
$ws.Range("G9").Value = -20.40935985752628
$ws.Range("H9").Value = -20.40935985752628
$ws.Range("I9").Value = -20.40935985752628
$ws.Range("J9").Value = -20.40935985752628
$ws.Range("K9").Value = -20.40935985752628
$ws.Range("B10").Value = -20.40935985752628
$ws.Range("C10").Value = -20.40935985752628
$ws.Range("D10").Value = -20.40935985752628
$ws.Range("E10").Value = -20.40935985752628
$ws.Range("F10").Value = -20.40935985752628
$ws.Range("G10").Value = -20.40935985752628
$ws.Range("H10").Value = -20.40935985752628
$ws.Range("I10").Value = 1.715999691000212
$ws.Range("J10").Value = -20.40935985752628
$ws.Range("K10").Value = 2.232447817232296
$ws.Range("B11").Value = -20.40935985752628
$ws.Range("C11").Value = -20.40935985752628
$ws.Range("D11").Value = -20.40935985752628
$ws.Range("E11").Value = 2.907024021714002
$ws.Range("F11").Value = -20.40935985752628
$ws.Range("G11").Value = 2.821207281462056
$ws.Range("H11").Value = -20.40935985752628
$ws.Range("I11").Value = -20.40935985752628
$ws.Range("J11").Value = -20.40935985752628
$ws.Range("K11").Value = 1.962789162053875
$ws.Range("B12").Value = -20.40935985752628
$ws.Range("C12").Value = -20.40935985752628
$ws.Range("D12").Value = -20.40935985752628
$ws.Range("E12").Value = -20.40935985752628
$ws.Range("F12").Value = -20.40935985752628
$ws.Range("G12").Value = -20.40935985752628
$ws.Range("H12").Value = -20.40935985752628
$ws.Range("I12").Value = -20.40935985752628
$ws.Range("J12").Value = -20.40935985752628
$ws.Range("K12").Value = -20.40935985752628
$ws.Range("B13").Value = -20.40935985752628
$ws.Range("C13").Value = -20.40935985752628
$ws.Range("D13").Value = -20.40935985752628
$ws.Range("E13").Value = 2.458332323733351
$ws.Range("F13").Value = -20.40935985752628
$ws.Range("G13").Value = -20.40935985752628
$ws.Range("H13").Value = -20.40935985752628
$ws.Range("I13").Value = -20.40935985752628
$ws.Range("J13").Value = 1.717506224422876
$ws.Range("K13").Value = 1.809036940088381
$ws.Range("B14").Value = -20.40935985752628
$ws.Range("C14").Value = -20.40935985752628
$ws.Range("D14").Value = 1.62477984216458
$ws.Range("E14").Value = -20.40935985752628
$ws.Range("F14").Value = -20.40935985752628
$ws.Range("G14").Value = -20.40935985752628
$ws.Range("H14").Value = -20.40935985752628
$ws.Range("I14").Value = -20.40935985752628
$ws.Range("J14").Value = -20.40935985752628
$ws.Range("K14").Value = 1.969105673640736
$ws.Range("B15").Value = -20.40935985752628
$ws.Range("C15").Value = -20.40935985752628
$ws.Range("D15").Value = 1.604264395517794
$ws.Range("E15").Value = -20.40935985752628
$ws.Range("F15").Value = -20.40935985752628
$ws.Range("G15").Value = -20.40935985752628
$ws.Range("H15").Value = -20.40935985752628
$ws.Range("I15").Value = -20.40935985752628
$ws.Range("J15").Value = -20.40935985752628
$ws.Range("K15").Value = -20.40935985752628
$ws.Range("B16").Value = -20.40935985752628
$ws.Range("C16").Value = -20.40935985752628
$ws.Range("D16").Value = -20.40935985752628
$ws.Range("E16").Value = -20.40935985752628
$ws.Range("F16").Value = -20.40935985752628
$ws.Range("G16").Value = -20.40935985752628
$ws.Range("H16").Value = -20.40935985752628
$ws.Range("I16").Value = -20.40935985752628
$ws.Range("J16").Value = 1.897470891967784
$ws.Range("K16").Value = -20.40935985752628
$ws.Range("B17").Value = -20.40935985752628
$ws.Range("C17").Value = 2.330120200671251
$ws.Range("D17").Value = 1.935131384056204
$ws.Range("E17").Value = -20.40935985752628
$ws.Range("F17").Value = -20.40935985752628
$ws.Range("G17").Value = -20.40935985752628
$ws.Range("H17").Value = 2.083212420906162
$ws.Range("I17").Value = 2.119868834908839
$ws.Range("J17").Value = 2.491114333322554
$ws.Range("K17").Value = -20.40935985752628
$ws.Range("B18").Value = -20.40935985752628
$ws.Range("C18").Value = -20.40935985752628
$ws.Range("D18").Value = -20.40935985752628
$ws.Range("E18").Value = -20.40935985752628
$ws.Range("F18").Value = -20.40935985752628
$ws.Range("G18").Value = -20.40935985752628
$ws.Range("H18").Value = 2.016429986918027
$ws.Range("I18").Value = 2.053717685323813
$ws.Range("J18").Value = 2.400647839700246
$ws.Range("K18").Value = -20.40935985752628
$ws.Range("B19").Value = -20.40935985752628
$ws.Range("C19").Value = -20.40935985752628
$ws.Range("D19").Value = 2.055617815458025
$ws.Range("E19").Value = -20.40935985752628
$ws.Range("F19").Value = -20.40935985752628
$ws.Range("G19").Value = -20.40935985752628
$ws.Range("H19").Value = 1.636889230639025
$ws.Range("I19").Value = 1.810490781070056
$ws.Range("J19").Value = -20.40935985752628
$ws.Range("K19").Value = -20.40935985752628
$ws.Range("B20").Value = -20.40935985752628
$ws.Range("C20").Value = 0.8785446807132508
$ws.Range("D20").Value = 1.412554616597077
$ws.Range("E20").Value = -20.40935985752628
$ws.Range("F20").Value = 4.321927110720398
$ws.Range("G20").Value = -20.40935985752628
$ws.Range("H20").Value = 1.623334656915258
$ws.Range("I20").Value = 1.192541931002583
$ws.Range("J20").Value = -20.40935985752628
$ws.Range("K20").Value = 1.993927100732957
$ws.Range("B21").Value = -20.40935985752628
$ws.Range("C21").Value = 1.191199636325147
$ws.Range("D21").Value = -20.40935985752628
$ws.Range("E21").Value = 1.85125220399174
$ws.Range("F21").Value = -20.40935985752628
$ws.Range("G21").Value = 2.447722277947184
$ws.Range("H21").Value = 1.449529109991784
$ws.Range("I21").Value = -20.40935985752628
$ws.Range("J21").Value = -20.40935985752628
$ws.Range("K21").Value = -20.40935985752628
